{"js": "// 1. Title paragraph \"\u4e00\u3001\u6211\u7684\u65e9\u5e74\u751f\u6d3b\": make it red and size 24pt (sz=48 half-points).\nconst titleResults = context.document.body.search(\"\u4e00\u3001\u6211\u7684\u65e9\u5e74\u751f\u6d3b\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nconst title = titleResults.items[0];\ntitle.font.color = \"#FF0000\";\ntitle.font.size = 24;\n\n// 2. \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4f9d\u8cf4\u65bc\u767c\u660e\u3002\" -> \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4ef0\u8cf4\u767c\u660e\u3002\"\nconst r1 = context.document.body.search(\"\u4f9d\u8cf4\u65bc\u767c\u660e\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nr1.items[0].insertText(\"\u4ef0\u8cf4\u767c\u660e\", Word.InsertLocation.replace);\n\n// 3. \"\u9019\u6b63\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\" -> \"\u9019\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\"\nconst r2 = context.document.body.search(\"\u9019\u6b63\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\", { matchCase: true });\nr2.load(\"items\");\nawait context.sync();\nr2.items[0].insertText(\"\u9019\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\", Word.InsertLocation.replace);\n\n// 4. \"\u5176\u6700\u7d42\u76ee\u7684\u662f\u5b8c\u5168\u638c\u63e1\u7cbe\u795e\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u63a7\u5236\uff0c\u99d5\u99ad\u81ea\u7136\u7684\u529b\u91cf\u4ee5\u6eff\u8db3\u4eba\u985e\u7684\u9700\u6c42\u3002\"\n//    -> \"\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\uff0c\u99d5\u99ad\u7269\u8cea\u4e16\u754c\uff0c\u5c07\u81ea\u7136\u4e4b\u529b\u904b\u7528\u65bc\u4eba\u985e\u9700\u6c42\u3002\"\nconst r3 = context.document.body.search(\"\u5176\u6700\u7d42\u76ee\u7684\u662f\u5b8c\u5168\u638c\u63e1\u7cbe\u795e\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u63a7\u5236\uff0c\u99d5\u99ad\u81ea\u7136\u7684\u529b\u91cf\u4ee5\u6eff\u8db3\u4eba\u985e\u7684\u9700\u6c42\u3002\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nr3.items[0].insertText(\"\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\uff0c\u99d5\u99ad\u7269\u8cea\u4e16\u754c\uff0c\u5c07\u81ea\u7136\u4e4b\u529b\u904b\u7528\u65bc\u4eba\u985e\u9700\u6c42\u3002\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Title paragraph \"\u4e00\u3001\u6211\u7684\u65e9\u5e74\u751f\u6d3b\": make it red and size 24pt (sz=48 half-points).\n$titleRange = $d.Content\n$titleFind = $titleRange.Find\n$titleFind.ClearFormatting()\n$titleFind.Text = \"\u4e00\u3001\u6211\u7684\u65e9\u5e74\u751f\u6d3b\"\n$titleFind.MatchCase = $true\n$titleFound = $titleFind.Execute()\nif ($titleFound) {\n    $titleRange.Font.Color = 255\n    $titleRange.Font.Size = 24\n}\n\n# 2. \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4f9d\u8cf4\u65bc\u767c\u660e\u3002\" -> \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4ef0\u8cf4\u767c\u660e\u3002\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"\u4f9d\u8cf4\u65bc\u767c\u660e\"\n$find1.Replacement.Text = \"\u4ef0\u8cf4\u767c\u660e\"\n$find1.MatchCase = $true\n[void]$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2)\n\n# 3. \"\u9019\u6b63\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\" -> \"\u9019\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"\u9019\u6b63\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\"\n$find2.Replacement.Text = \"\u9019\u662f\u4ed6\u5275\u610f\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\"\n$find2.MatchCase = $true\n[void]$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n# 4. \"\u5176\u6700\u7d42\u76ee\u7684\u662f\u5b8c\u5168\u638c\u63e1\u7cbe\u795e\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u63a7\u5236\uff0c\u99d5\u99ad\u81ea\u7136\u7684\u529b\u91cf\u4ee5\u6eff\u8db3\u4eba\u985e\u7684\u9700\u6c42\u3002\"\n#    -> \"\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\uff0c\u99d5\u99ad\u7269\u8cea\u4e16\u754c\uff0c\u5c07\u81ea\u7136\u4e4b\u529b\u904b\u7528\u65bc\u4eba\u985e\u9700\u6c42\u3002\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"\u5176\u6700\u7d42\u76ee\u7684\u662f\u5b8c\u5168\u638c\u63e1\u7cbe\u795e\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u63a7\u5236\uff0c\u99d5\u99ad\u81ea\u7136\u7684\u529b\u91cf\u4ee5\u6eff\u8db3\u4eba\u985e\u7684\u9700\u6c42\u3002\"\n$find3.Replacement.Text = \"\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\uff0c\u99d5\u99ad\u7269\u8cea\u4e16\u754c\uff0c\u5c07\u81ea\u7136\u4e4b\u529b\u904b\u7528\u65bc\u4eba\u985e\u9700\u6c42\u3002\"\n$find3.MatchCase = $true\n[void]$find3.Execute([ref]$find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find3.Replacement.Text, 2)\n"}
